$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 108, shifting existing rows 108+ down by one.
$ws.Rows.Item(108).Insert()

# Populate the newly inserted row 108 with the new data record.
$ws.Cells.Item(108, 1).Value = 8
$ws.Cells.Item(108, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(108, 3).Value = "Coquimbo"
$ws.Cells.Item(108, 4).Value = 44566
$ws.Cells.Item(108, 5).Value = 4
$ws.Cells.Item(108, 6).Value = 100112012
$ws.Cells.Item(108, 7).Value = "Espinaca"
$ws.Cells.Item(108, 8).Value = "Sin especificar"
$ws.Cells.Item(108, 9).Value = "Primera"
$ws.Cells.Item(108, 10).Value = 2960
$ws.Cells.Item(108, 11).Value = 400
$ws.Cells.Item(108, 12).Value = 500
$ws.Cells.Item(108, 13).Value = 450
$ws.Cells.Item(108, 14).Value = '$/atado 300 a 500 gramos'
$ws.Cells.Item(108, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(108, 16).Value = 900
$ws.Cells.Item(108, 17).Value = 0.5
$ws.Cells.Item(108, 18).Value = "Hortaliza"
